# Saldo_guide.xlsx update:
#  - Sheet renamed from IClientBalance-20240618-094409- to IClientBalance-20240619-111710-
#  - "Dt. Referencia" (column G) bumped one day, 45461 -> 45462, for every data row (2..257)
#  - Row 57 (PHILIPE FERREIRA DA SILVA LIMA): Saldo Previsto (D) / Vl. Total (H) corrected
#    from 19.21 to 3013.33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 257

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45462
}

$ws.Cells.Item(57, 4).Value = 3013.33
$ws.Cells.Item(57, 8).Value = 3013.33

$ws.Name = "IClientBalance-20240619-111710-"
